$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @{
    2 = @(21.70413948077697, 8.942990563587324, 6.010979469025904, 10.33024867335824, 0, 49.4585704163182, 18.84614408008562, 0, 0, 0, 10.21525176535112, 0, 19.01868609629499)
    3 = @(21.19981760106075, 8.451239213540052, 5.898105816394818, 10.35533621503886, 0, 48.93872331405137, 18.83083664099944, 0, 0, 0, 10.19624004958483, 0, 19.09193985957641)
    4 = @(20.89057683872107, 8.132762874485469, 5.829692564669692, 10.371850108656, 0, 48.63482696125696, 18.82624201564692, 0, 0, 0, 10.18671794623453, 0, 19.13891830692182)
    5 = @(20.76486293406493, 7.998870930166772, 5.802078355662226, 10.3788590910741, 0, 48.51495817835218, 18.825575891681, 0, 0, 0, 10.18338057971497, 0, 19.15856669350289)
    6 = @(20.74401254690908, 7.976391680080673, 5.797510209459879, 10.38003981369918, 0, 48.49529741344617, 18.82553806395159, 0, 0, 0, 10.18285925735103, 0, 19.16185978878325)
    7 = @(20.88887991761187, 8.130973730199505, 5.829319024328815, 10.37194350239906, 0, 48.63319412969016, 18.82622815135569, 0, 0, 0, 10.1866707365286, 0, 19.13918124846413)
    8 = @(21.53028143070796, 8.776893045192827, 5.971899903492547, 10.33866865300502, 0, 49.27623201435844, 18.8398681219857, 0, 0, 0, 10.20825140373828, 0, 19.04352960101993)
    9 = @(22.78288537457854, 9.910715341230295, 6.256801044148781, 10.28221078205235, 0, 50.65255694349319, 18.90477979724658, 0, 0, 0, 10.267531765165, 0, 18.87176858742344)
    10 = @(23.6892211973355, 10.66124505360214, 6.467014322534332, 10.24607364442369, 0, 51.72545965843448, 18.97574275818932, 0, 0, 0, 10.3212541368396, 0, 18.75512849128397)
    11 = @(24.09641265778603, 10.98458586009791, 6.562366549564433, 10.23079017531595, 0, 52.22500343938806, 19.01306098785128, 0, 0, 0, 10.3478565845428, 0, 18.70412186600608)
    12 = @(24.24971142741044, 11.10441768893834, 6.598396939674153, 10.22516861074876, 0, 52.41564921077446, 19.02791350800817, 0, 0, 0, 10.35823684015091, 0, 18.68510098214932)
    13 = @(24.21673808771013, 11.07872598057467, 6.590641178487046, 10.2263719376311, 0, 52.37452717118297, 19.0246827529416, 0, 0, 0, 10.35598771415822, 0, 18.68918440128353)
    14 = @(24.10904347828717, 10.99449685756336, 6.565332553936493, 10.230324361011, 0, 52.24065931611366, 19.01426848811511, 0, 0, 0, 10.34870446022841, 0, 18.70255111597289)
    15 = @(24.04295615809648, 10.94256390319393, 6.549819112788942, 10.23276694052599, 0, 52.15884895448009, 19.00798321833116, 0, 0, 0, 10.34428302460916, 0, 18.71077689894366)
    16 = @(23.66249346656945, 10.639749176038, 6.460774018835039, 10.24709568888182, 0, 51.69302965320028, 18.97340501151761, 0, 0, 0, 10.31955875669456, 0, 18.75850312395312)
    17 = @(23.4276657122423, 10.44934280749293, 6.40604981690818, 10.25618171366792, 0, 51.41008279786979, 18.95348035875512, 0, 0, 0, 10.30494224384155, 0, 18.78830680319042)
    18 = @(23.29212742366401, 10.33812476246123, 6.374550195175688, 10.26151652742318, 0, 51.24842971220225, 18.94249457860249, 0, 0, 0, 10.2967391928397, 0, 18.80564244317361)
    19 = @(23.2461606072259, 10.30017642929794, 6.363882076456822, 10.2633414918667, 0, 51.19388921210864, 18.93885653971903, 0, 0, 0, 10.29399695022937, 0, 18.81154523457421)
    20 = @(23.45271354122393, 10.46978806236159, 6.411878035575935, 10.25520323478026, 0, 51.44009134586955, 18.95555229970158, 0, 0, 0, 10.30647711793095, 0, 18.78511414735666)
    21 = @(24.14070153899809, 11.01930790015174, 6.572768705002696, 10.22915893612781, 0, 52.27994076839691, 19.01730787629317, 0, 0, 0, 10.35083544981963, 0, 18.69861701133201)
    22 = @(24.58505319124288, 11.3632337786388, 6.677452468391765, 10.21310469514596, 0, 52.83737388940977, 19.06186853117067, 0, 0, 0, 10.3816101634206, 0, 18.64380057743704)
    23 = @(24.34842867241522, 11.18106890536639, 6.621635607438442, 10.22158471155416, 0, 52.53913639737168, 19.03770274846256, 0, 0, 0, 10.36502353555771, 0, 18.67290061818546)
    24 = @(23.44139106454349, 10.46055021541295, 6.409243213223374, 10.25564525890257, 0, 51.42652130119176, 18.95461411360147, 0, 0, 0, 10.30578257749892, 0, 18.78655692057148)
    25 = @(22.44572724524829, 9.618391066015537, 6.179404917834677, 10.29654462269822, 0, 50.26873905238733, 18.88312931549296, 0, 0, 0, 10.24969494174331, 0, 18.91655076700902)
}

foreach ($rowNum in $rowsData.Keys) {
    $rowValues = $rowsData[$rowNum]
    $arr = New-Object "object[,]" 1, $rowValues.Length
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $arr[0, $i] = $rowValues[$i]
    }
    $rangeAddr = "B" + $rowNum + ":N" + $rowNum
    $ws.Range($rangeAddr).Value = $arr
}